$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert new "Index_Benchmark" column at A, shift Role_Benchmark to B,
# drop old Overall_Benchmark column (values for C:G are columns D:H pre-edit).
$header = New-Object 'object[,]' 1,7
$header[0,0] = "Index_Benchmark"
$header[0,1] = "Role_Benchmark"
$header[0,2] = "Work_Life_Balance_Benchmark"
$header[0,3] = "Career_Growth_Benchmark"
$header[0,4] = "Compensation_Benchmark"
$header[0,5] = "Leadership_Benchmark"
$header[0,6] = "Colleagues_Benchmark"
$ws.Range("A1:G1").Value = $header

$data = New-Object 'object[,]' 14,7
$data[0,0] = 1
$data[0,1] = "Accountant"
$data[0,2] = 42.66457320065805
$data[0,3] = 83.44394711474659
$data[0,4] = 110.3736731718898
$data[0,5] = 26.83152524691578
$data[0,6] = 104.8291539435402

$data[1,0] = 2
$data[1,1] = "Business Analyst"
$data[1,2] = 32.07033993460573
$data[1,3] = 51.42882730440155
$data[1,4] = 27.10409715298731
$data[1,5] = 41.3234115956139
$data[1,6] = 81.89001172060458

$data[2,0] = 3
$data[2,1] = "Customer Support"
$data[2,2] = 85.90819579932058
$data[2,3] = 91.30494565446023
$data[2,4] = 114.3055978810442
$data[2,5] = 75.29002474858609
$data[2,6] = 91.59631169145987

$data[3,0] = 4
$data[3,1] = "Data Analyst"
$data[3,2] = 85.24068197080739
$data[3,3] = 53.85237101017578
$data[3,4] = 88.08546543371043
$data[3,5] = 39.26178870883071
$data[3,6] = 57.43328769992886

$data[4,0] = 5
$data[4,1] = "Data Scientist"
$data[4,2] = 58.29060015310971
$data[4,3] = 24.26173277504419
$data[4,4] = 62.93630877370196
$data[4,5] = 104.9947959832414
$data[4,6] = 80.61764733720652

$data[5,0] = 6
$data[5,1] = "HR Specialist"
$data[5,2] = 98.02336050081948
$data[5,3] = 104.2198591589407
$data[5,4] = 47.63064210578121
$data[5,5] = 38.5117707495847
$data[5,6] = 123.9754787741175

$data[6,0] = 7
$data[6,1] = "Legal Advisor"
$data[6,2] = 41.49928456059973
$data[6,3] = 39.41613945827162
$data[6,4] = 117.4910927359335
$data[6,5] = 101.4742991576342
$data[6,6] = 68.95778873870678

$data[7,0] = 8
$data[7,1] = "Machine Learning Engineer"
$data[7,2] = 55.33196789161902
$data[7,3] = 108.0914655182148
$data[7,4] = 121.3040254051418
$data[7,5] = 106.6021677754157
$data[7,6] = 33.11104989979082

$data[8,0] = 9
$data[8,1] = "Marketing Manager"
$data[8,2] = 76.93306421135722
$data[8,3] = 11.83680526988351
$data[8,4] = 46.35271529981099
$data[8,5] = 69.34155507426746
$data[8,6] = 59.86941548676152

$data[9,0] = 10
$data[9,1] = "Operations Manager"
$data[9,2] = 117.8684676525102
$data[9,3] = 34.93540751392739
$data[9,4] = 66.26929259530763
$data[9,5] = 84.74013383104247
$data[9,6] = 66.28163934977103

$data[10,0] = 11
$data[10,1] = "Product Designer"
$data[10,2] = 26.7311540241746
$data[10,3] = 58.33904252544976
$data[10,4] = 118.9597437623324
$data[10,5] = 108.5510808987087
$data[10,6] = 44.66260834377468

$data[11,0] = 12
$data[11,1] = "Project Manager"
$data[11,2] = 57.81059284620971
$data[11,3] = 28.28673906268114
$data[11,4] = 120.7498675189967
$data[11,5] = 27.79284006792783
$data[11,6] = 47.0029069858396

$data[12,0] = 13
$data[12,1] = "Sales Executive"
$data[12,2] = 115.7071936317938
$data[12,3] = 67.684081817014
$data[12,4] = 76.24929741873757
$data[12,5] = 73.53801272557061
$data[12,6] = 51.17971753894523

$data[13,0] = 14
$data[13,1] = "Software Engineer"
$data[13,2] = 90.0491947810874
$data[13,3] = 47.67084159234977
$data[13,4] = 45.79708627600634
$data[13,5] = 29.98705293764506
$data[13,6] = 78.56641314753028

$ws.Range("A2:G15").Value = $data
